# Applies the BetExplorer "liga-portugal-2 2023-2024" update:
#  1. Swap the two Nov-04 matches currently stored in rows 74 and 75
#     (Mafra-Leixoes <-> FC Porto B-Feirense), columns F:V only.
#  2. Append two new matches as rows 88 and 89 (Torreense-Mafra,
#     Vilaverdense-Tondela).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Swap rows 74 and 75 (columns F..V) -- A..E (index/meta/date) stay put
# ---------------------------------------------------------------------

$row74 = @("Mafra", 0, "Leixoes", 1, 1.88, "01/11/2023 16:12", 1.93, "04/11/2023 11:48", 3.73, "01/11/2023 16:12", 3.55, "04/11/2023 11:51", 3.8, "01/11/2023 16:12", 4.21, "04/11/2023 11:51", "https://www.betexplorer.com/football/portugal/liga-portugal-2/mafra-leixoes/YiBBPnTT/")
$row75 = @("FC Porto B", 2, "Feirense", 0, 1.98, "01/11/2023 16:12", 1.84, "04/11/2023 11:59", 3.57, "01/11/2023 16:12", 3.78, "04/11/2023 11:59", 3.87, "01/11/2023 16:12", 4.41, "04/11/2023 11:58", "https://www.betexplorer.com/football/portugal/liga-portugal-2/fc-porto-feirense/jTL6QSDN/")

$newRow74 = $row75
$newRow75 = $row74

for ($i = 0; $i -lt $newRow74.Length; $i++) {
    $col = 6 + $i
    $ws.Cells.Item(74, $col).Value = $newRow74[$i]
    $ws.Cells.Item(75, $col).Value = $newRow75[$i]
}

# ---------------------------------------------------------------------
# 2) Append rows 88 and 89
# ---------------------------------------------------------------------

$ws.Cells.Item(88, 1).Value = 87
$ws.Cells.Item(88, 2).Value = "portugal"
$ws.Cells.Item(88, 3).Value = "liga-portugal-2"
$ws.Cells.Item(88, 4).Value = "2023-2024"
$ws.Cells.Item(88, 5).Value = 45242.5
$ws.Cells.Item(88, 6).Value = "Torreense"
$ws.Cells.Item(88, 7).Value = 0
$ws.Cells.Item(88, 8).Value = "Mafra"
$ws.Cells.Item(88, 9).Value = 0
$ws.Cells.Item(88, 10).Value = 2.12
$ws.Cells.Item(88, 11).Value = "08/11/2023 06:12"
$ws.Cells.Item(88, 12).Value = 2.57
$ws.Cells.Item(88, 13).Value = "12/11/2023 11:50"
$ws.Cells.Item(88, 14).Value = 3.41
$ws.Cells.Item(88, 15).Value = "08/11/2023 06:12"
$ws.Cells.Item(88, 16).Value = 3.18
$ws.Cells.Item(88, 17).Value = "12/11/2023 11:50"
$ws.Cells.Item(88, 18).Value = 3.58
$ws.Cells.Item(88, 19).Value = "08/11/2023 06:12"
$ws.Cells.Item(88, 20).Value = 3.04
$ws.Cells.Item(88, 21).Value = "12/11/2023 11:50"
$ws.Cells.Item(88, 22).Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/torreense-mafra/jTsxt8DA/"

$ws.Cells.Item(89, 1).Value = 88
$ws.Cells.Item(89, 2).Value = "portugal"
$ws.Cells.Item(89, 3).Value = "liga-portugal-2"
$ws.Cells.Item(89, 4).Value = "2023-2024"
$ws.Cells.Item(89, 5).Value = 45242.5
$ws.Cells.Item(89, 6).Value = "Vilaverdense"
$ws.Cells.Item(89, 7).Value = 1
$ws.Cells.Item(89, 8).Value = "Tondela"
$ws.Cells.Item(89, 9).Value = 2
$ws.Cells.Item(89, 10).Value = 2.68
$ws.Cells.Item(89, 11).Value = "08/11/2023 06:12"
$ws.Cells.Item(89, 12).Value = 3.73
$ws.Cells.Item(89, 13).Value = "12/11/2023 11:57"
$ws.Cells.Item(89, 14).Value = 3.34
$ws.Cells.Item(89, 15).Value = "08/11/2023 06:12"
$ws.Cells.Item(89, 16).Value = 3.49
$ws.Cells.Item(89, 17).Value = "12/11/2023 11:57"
$ws.Cells.Item(89, 18).Value = 2.58
$ws.Cells.Item(89, 19).Value = "08/11/2023 06:12"
$ws.Cells.Item(89, 20).Value = 2.09
$ws.Cells.Item(89, 21).Value = "12/11/2023 11:57"
$ws.Cells.Item(89, 22).Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/vilaverdense-fc-tondela/UJttuSSG/"

# Match style of existing data rows: bold/bordered index column (A) and
# the date-time numeric format on column E (format-only copy so the
# values just written above are preserved).
$ws.Range("A87").Copy() | Out-Null
$ws.Range("A88:A89").PasteSpecial(-4122) | Out-Null
$ws.Range("E87").Copy() | Out-Null
$ws.Range("E88:E89").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3) Update the used-range dimension
# ---------------------------------------------------------------------
$ws.UsedRange | Out-Null
